$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urlPN2222   = "https://www.digikey.com/product-detail/en/on-semiconductor/PN2222ATFR/PN2222AD26ZCT-ND/459004"
$urlRHM10K   = "https://www.digikey.com/product-detail/en/rohm-semiconductor/ESR03EZPJ103/RHM10KDCT-ND/1762925"
$urlRHM1     = "https://www.digikey.com/product-detail/en/rohm-semiconductor/ESR03EZPJ1R0/RHM1DCT-ND/4053738"
$urlTDK      = "https://www.digikey.com/product-detail/en/tdk-corporation/C1608X5R1V475K080AC/445-9064-1-ND/3648692"
$urlTDKspace = "https://www.digikey.com/product-detail/en/tdk-corporation/C1608X5R1V475K080AC/445-9064-1-ND/3648692 "
$urlNSR      = "https://www.digikey.com/product-detail/en/on-semiconductor/NSR05T40P2T5G/NSR05T40P2T5GOSCT-ND/5761655"

# --- Phase 1: set cell values in the exact order that reproduces the
# sharedStrings.xml append order from the target diff ---
$ws.Range("C10").Value = "PN2222AD26ZCT-ND"
$ws.Range("B10").Value = "PN2222ATFR"
$ws.Range("D10").Value = $urlPN2222
$ws.Range("A10").Value = "npn transistor "

$ws.Range("C12").Value = "RHM1DCT-ND "
$ws.Range("A12").Value = "1ohm 0603 "
$ws.Range("B12").Value = "ESR03EZPJ1R0"

$ws.Range("C11").Value = "RHM10KDCT-ND "
$ws.Range("A11").Value = "10kohm 0603"
$ws.Range("B11").Value = "ESR03EZPJ103"

$ws.Range("D11").Value = $urlRHM10K
$ws.Range("D12").Value = $urlRHM1

$ws.Range("A13").Value = "CAP CER 4.7UF 35V X5R 0603"
$ws.Range("C13").Value = "445-9064-1-ND "
$ws.Range("B13").Value = "C1608X5R1V475K080AC"
$ws.Range("D13").Value = $urlTDKspace

# --- Phase 2: apply the Hyperlink cell style to all the new cells that need it ---
$styled = @("A10","B10","D10","A11","B11","C11","D11","A12","B12","C12","D12","A13","B13","C13","D13","D3")
foreach ($addr in $styled) {
    $ws.Range($addr).Style = "Hyperlink"
}

# --- Phase 3: add the hyperlinks themselves, in the order that reproduces
# the target rId numbering. For cells whose hyperlink carries a "display"
# attribute, Excel's Hyperlinks.Add also overwrites the cell text with
# TextToDisplay, so immediately restore the cell's real value afterwards. ---

$ws.Hyperlinks.Add($ws.Range("B10"), $urlPN2222, "", "", $urlPN2222)
$ws.Range("B10").Value = "PN2222ATFR"

$ws.Hyperlinks.Add($ws.Range("A10"), $urlPN2222, "", "", $urlPN2222)
$ws.Range("A10").Value = "npn transistor "

$ws.Hyperlinks.Add($ws.Range("A12"), $urlRHM1, "", "", $urlRHM1)
$ws.Range("A12").Value = "1ohm 0603 "

$ws.Hyperlinks.Add($ws.Range("C12"), $urlRHM1, "", "", $urlRHM1)
$ws.Range("C12").Value = "RHM1DCT-ND "

$ws.Hyperlinks.Add($ws.Range("B12"), $urlRHM1, "", "", $urlRHM1)
$ws.Range("B12").Value = "ESR03EZPJ1R0"

$ws.Hyperlinks.Add($ws.Range("A11"), $urlRHM10K, "", "", $urlRHM10K)
$ws.Range("A11").Value = "10kohm 0603"

$ws.Hyperlinks.Add($ws.Range("C11"), $urlRHM10K, "", "", $urlRHM10K)
$ws.Range("C11").Value = "RHM10KDCT-ND "

$ws.Hyperlinks.Add($ws.Range("B11"), $urlRHM10K, "", "", $urlRHM10K)
$ws.Range("B11").Value = "ESR03EZPJ103"

$ws.Hyperlinks.Add($ws.Range("A13"), $urlTDK, "", "", $urlTDK)
$ws.Range("A13").Value = "CAP CER 4.7UF 35V X5R 0603"

$ws.Hyperlinks.Add($ws.Range("C13"), $urlTDK, "", "", $urlTDK)
$ws.Range("C13").Value = "445-9064-1-ND "

$ws.Hyperlinks.Add($ws.Range("B13"), $urlTDK, "", "", $urlTDK)
$ws.Range("B13").Value = "C1608X5R1V475K080AC"

# These four carry no display text - the cell keeps its own value (the URL).
$ws.Hyperlinks.Add($ws.Range("D12"), $urlRHM1)
$ws.Hyperlinks.Add($ws.Range("D11"), $urlRHM10K)
$ws.Hyperlinks.Add($ws.Range("D10"), $urlPN2222)
$ws.Hyperlinks.Add($ws.Range("D13"), $urlTDK)

$ws.Hyperlinks.Add($ws.Range("D3"), $urlNSR)

# Re-apply the Hyperlink cell style, since Hyperlinks.Add touches formatting.
foreach ($addr in $styled) {
    $ws.Range($addr).Style = "Hyperlink"
}

# --- Phase 4: window / selection state ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
[void]$ws.Range("D21").Select()

Write-Host "done"
